# Auto-generated Excel COM-interop script to apply the Chocobo_Profits.xlsx numeric updates
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value2 = 234.4
$ws.Range("I33").Value2 = 208.14285
$ws.Range("J33").Value2 = 602
$ws.Range("K33").Value2 = 208.14285
$ws.Range("L33").Value2 = 602
$ws.Range("M33").Value2 = 20.85714999999999
$ws.Range("N33").Value2 = -1060

$ws.Range("H129").Value2 = 988.3889
$ws.Range("J129").Value2 = 988.3889
$ws.Range("L129").Value2 = 2965.1667
$ws.Range("N129").Value2 = -12965.1667

$ws.Range("H137").Value2 = 735411.4
$ws.Range("I137").Value2 = 2382485.5
$ws.Range("J137").Value2 = 3378.4888
$ws.Range("K137").Value2 = 7147456.5
$ws.Range("L137").Value2 = 10135.4664
$ws.Range("M137").Value2 = -7144906.5
$ws.Range("N137").Value2 = -15235.4664

$ws.Range("H138").Value2 = 3380.5571
$ws.Range("I138").Value2 = 2472.9443
$ws.Range("J138").Value2 = 3694.7307
$ws.Range("K138").Value2 = 7418.8329
$ws.Range("L138").Value2 = 11084.1921
$ws.Range("M138").Value2 = -2278.8329
$ws.Range("N138").Value2 = -21364.1921

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value2 = 4065.2527
$ws.Range("I32").Value2 = 2819.178
$ws.Range("K32").Value2 = 2819.178
$ws.Range("M32").Value2 = -2532.178

$ws.Range("H61").Value2 = 1304.2
$ws.Range("I61").Value2 = 1304.6666
$ws.Range("J61").Value2 = 1300
$ws.Range("K61").Value2 = 1304.6666
$ws.Range("L61").Value2 = 1300
$ws.Range("M61").Value2 = -1092.6666
$ws.Range("N61").Value2 = -1724

$ws.Range("H74").Value2 = 1287.6182
$ws.Range("I74").Value2 = 798.61365
$ws.Range("J74").Value2 = 3243.6365
$ws.Range("K74").Value2 = 798.61365
$ws.Range("L74").Value2 = 3243.6365
$ws.Range("M74").Value2 = 75.38634999999999
$ws.Range("N74").Value2 = -4991.636500000001

$ws.Range("H77").Value2 = 1287.6182
$ws.Range("I77").Value2 = 798.61365
$ws.Range("J77").Value2 = 3243.6365
$ws.Range("K77").Value2 = 3993.06825
$ws.Range("L77").Value2 = 16218.1825
$ws.Range("M77").Value2 = 374.9317499999997
$ws.Range("N77").Value2 = -24954.1825

$ws.Range("H136").Value2 = 1304.2
$ws.Range("I136").Value2 = 1304.6666
$ws.Range("J136").Value2 = 1300
$ws.Range("K136").Value2 = 3913.9998
$ws.Range("L136").Value2 = 3900
$ws.Range("M136").Value2 = -1363.9998
$ws.Range("N136").Value2 = -9000

$ws.Range("H138").Value2 = 0
$ws.Range("J138").Value2 = 0
$ws.Range("L138").Value2 = 0
$ws.Range("N138").ClearContents()

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H64").Value2 = 521.8461
$ws.Range("I64").Value2 = 554
$ws.Range("J64").Value2 = 501.75
$ws.Range("K64").Value2 = 554
$ws.Range("L64").Value2 = 501.75
$ws.Range("M64").Value2 = -329
$ws.Range("N64").Value2 = -951.75

$ws.Range("H67").Value2 = 521.8461
$ws.Range("I67").Value2 = 554
$ws.Range("J67").Value2 = 501.75
$ws.Range("K67").Value2 = 554
$ws.Range("L67").Value2 = 501.75
$ws.Range("M67").Value2 = 226
$ws.Range("N67").Value2 = -2061.75

$ws.Range("H134").Value2 = 3512.1428
$ws.Range("I134").Value2 = 1781.8
$ws.Range("K134").Value2 = 5345.4
$ws.Range("M134").Value2 = -2810.4

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value2 = 3574.716
$ws.Range("I31").Value2 = 1485.7059
$ws.Range("J31").Value2 = 4129.6094
$ws.Range("K31").Value2 = 1485.7059
$ws.Range("L31").Value2 = 4129.6094
$ws.Range("M31").Value2 = -1190.7059
$ws.Range("N31").Value2 = -4719.6094

$ws.Range("H34").Value2 = 3574.716
$ws.Range("I34").Value2 = 1485.7059
$ws.Range("J34").Value2 = 4129.6094
$ws.Range("K34").Value2 = 1485.7059
$ws.Range("L34").Value2 = 4129.6094
$ws.Range("M34").Value2 = -1283.7059
$ws.Range("N34").Value2 = -4533.6094

$ws.Range("H132").Value2 = 2943.125
$ws.Range("I132").Value2 = 2576.1516
$ws.Range("J132").Value2 = 4673.143
$ws.Range("K132").Value2 = 7728.4548
$ws.Range("L132").Value2 = 14019.429
$ws.Range("M132").Value2 = -5198.4548
$ws.Range("N132").Value2 = -19079.429

$ws.Range("H134").Value2 = 6461.913
$ws.Range("I134").Value2 = 6983.8237
$ws.Range("J134").Value2 = 4983.1665
$ws.Range("K134").Value2 = 20951.4711
$ws.Range("L134").Value2 = 14949.4995
$ws.Range("M134").Value2 = -18416.4711
$ws.Range("N134").Value2 = -20019.4995

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H68").Value2 = 2165.6863
$ws.Range("I68").Value2 = 1255.2
$ws.Range("J68").Value2 = 3041.1538
$ws.Range("K68").Value2 = 3765.6
$ws.Range("L68").Value2 = 9123.4614
$ws.Range("M68").Value2 = -2954.6
$ws.Range("N68").Value2 = -10745.4614

$ws.Range("H71").Value2 = 2165.6863
$ws.Range("I71").Value2 = 1255.2
$ws.Range("J71").Value2 = 3041.1538
$ws.Range("K71").Value2 = 11296.8
$ws.Range("L71").Value2 = 27370.3842
$ws.Range("M71").Value2 = -7240.800000000001
$ws.Range("N71").Value2 = -35482.3842

$ws.Range("H107").Value2 = 16701844
$ws.Range("J107").Value2 = 41753616
$ws.Range("L107").Value2 = 125260848
$ws.Range("N107").Value2 = -125264688

$ws.Range("H113").Value2 = 3205679.2
$ws.Range("I113").Value2 = 554.2
$ws.Range("J113").Value2 = 8929117
$ws.Range("K113").Value2 = 1662.6
$ws.Range("L113").Value2 = 26787351
$ws.Range("M113").Value2 = 507.3999999999999
$ws.Range("N113").Value2 = -26791691

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H107").Value2 = 6536549
$ws.Range("I107").Value2 = 383.27274
$ws.Range("J107").Value2 = 18519520
$ws.Range("K107").Value2 = 383.27274
$ws.Range("L107").Value2 = 18519520
$ws.Range("M107").Value2 = 1536.72726
$ws.Range("N107").Value2 = -18523360

$ws.Range("H122").Value2 = 3221.55
$ws.Range("I122").Value2 = 2437.1177
$ws.Range("K122").Value2 = 7311.353099999999
$ws.Range("M122").Value2 = -4861.353099999999

$ws.Range("H126").Value2 = 3497.96
$ws.Range("I126").Value2 = 2854.2285
$ws.Range("J126").Value2 = 5000
$ws.Range("K126").Value2 = 8562.6855
$ws.Range("L126").Value2 = 15000
$ws.Range("M126").Value2 = -6092.6855
$ws.Range("N126").Value2 = -19940

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H122").Value2 = 4990.125
$ws.Range("I122").Value2 = 4360.3335
$ws.Range("J122").Value2 = 6879.5
$ws.Range("K122").Value2 = 13081.0005
$ws.Range("L122").Value2 = 20638.5
$ws.Range("M122").Value2 = -10631.0005
$ws.Range("N122").Value2 = -25538.5

$ws.Range("H132").Value2 = 4756.56
$ws.Range("I132").Value2 = 4040.7334
$ws.Range("J132").Value2 = 5830.3
$ws.Range("K132").Value2 = 12122.2002
$ws.Range("L132").Value2 = 17490.9
$ws.Range("M132").Value2 = -9592.200199999999
$ws.Range("N132").Value2 = -22550.9

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H98").Value2 = 50590
$ws.Range("J98").Value2 = 50590
$ws.Range("L98").Value2 = 50590
$ws.Range("N98").Value2 = -56580

$ws.Range("H122").Value2 = 3849.3794
$ws.Range("I122").Value2 = 2095.4
$ws.Range("J122").Value2 = 5728.643
$ws.Range("K122").Value2 = 6286.200000000001
$ws.Range("L122").Value2 = 17185.929
$ws.Range("M122").Value2 = -3836.200000000001
$ws.Range("N122").Value2 = -22085.929

$ws.Range("H132").Value2 = 11113671
$ws.Range("I132").Value2 = 1971.2778
$ws.Range("J132").Value2 = 27781220
$ws.Range("K132").Value2 = 5913.8334
$ws.Range("L132").Value2 = 83343660
$ws.Range("M132").Value2 = -3383.8334
$ws.Range("N132").Value2 = -83348720

Write-Host "Applied 182 cell updates."